# Refresh cryptocurrency Price / Volume(1h) columns (and the two coin swaps)
# from the latest coinranking.com scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.642.47"
$ws.Range("E2").Value = "  -2.38%  "

$ws.Range("D3").Value = "'3.480.71"
$ws.Range("E3").Value = "  -3.80%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'576.91"
$ws.Range("E5").Value = "  -4.23%  "

$ws.Range("D6").Value = "'188.60"
$ws.Range("E6").Value = "  -3.56%  "

$ws.Range("D7").Value = "'3.463.36"
$ws.Range("E7").Value = "  -3.95%  "

$ws.Range("D8").Value = "'0.604"
$ws.Range("E8").Value = "  -3.72%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("E10").Value = "  -5.64%  "

$ws.Range("E11").Value = "  -5.14%  "

$ws.Range("D12").Value = "'51.63"
$ws.Range("E12").Value = "  -3.05%  "

$ws.Range("E13").Value = "  -6.91%  "

$ws.Range("D14").Value = "'9.06"
$ws.Range("E14").Value = "  -5.40%  "

$ws.Range("D15").Value = "'4.032.57"
$ws.Range("E15").Value = "  -3.70%  "

$ws.Range("D16").Value = "'632.85"
$ws.Range("E16").Value = "  +5.07%  "

$ws.Range("D17").Value = "'68.677.98"
$ws.Range("E17").Value = "  -2.48%  "

# Row 18: coin identity changed (rank reshuffle)
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.487.88"
$ws.Range("E18").Value = "  -3.29%  "

# Row 19: coin identity changed (rank reshuffle)
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'12.40"
$ws.Range("E19").Value = "  -4.47%  "

$ws.Range("E20").Value = "  -2.35%  "

$ws.Range("D21").Value = "'18.06"
$ws.Range("E21").Value = "  -5.28%  "

$ws.Range("D22").Value = "'0.941"
$ws.Range("E22").Value = "  -5.86%  "

$ws.Range("D23").Value = "'17.77"
$ws.Range("E23").Value = "  -4.65%  "

$ws.Range("D24").Value = "'5.37"
$ws.Range("E24").Value = "  +3.06%  "

$ws.Range("D25").Value = "'99.09"
$ws.Range("E25").Value = "  -3.84%  "

$ws.Range("D26").Value = "'4.29"
$ws.Range("E26").Value = "  -7.05%  "

$ws.Range("E27").Value = "  -4.77%  "

$ws.Range("D29").Value = "'10.01"

$ws.Range("D30").Value = "'9.18"
$ws.Range("E30").Value = "  -5.60%  "

$ws.Range("D31").Value = "'32.41"
$ws.Range("E31").Value = "  -4.13%  "

$ws.Range("E32").Value = "  -8.42%  "

$ws.Range("D33").Value = "'4.04"
$ws.Range("E33").Value = "  -14.51%  "

$ws.Range("D34").Value = "'11.56"
$ws.Range("E34").Value = "  -5.94%  "

# Row 35: coin identity changed (rank reshuffle)
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.108"
$ws.Range("E35").Value = "  -7.77%  "

# Row 36: coin identity changed (rank reshuffle)
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'60.67"
$ws.Range("E36").Value = "  -4.17%  "

$ws.Range("D37").Value = "'3.701.31"
$ws.Range("E37").Value = "  -5.96%  "

$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.06%  "

$ws.Range("D39").Value = "'0.0₃0783"
$ws.Range("E39").Value = "  -11.45%  "

$ws.Range("D40").Value = "'499.82"
$ws.Range("E40").Value = "  -6.09%  "

$ws.Range("D41").Value = "'3.53"
$ws.Range("E41").Value = "  +0.10%  "

$ws.Range("D42").Value = "'2.91"
$ws.Range("E42").Value = "  -4.63%  "

$ws.Range("E43").Value = "  -5.77%  "

$ws.Range("E44").Value = "  -1.66%  "

$ws.Range("D45").Value = "'34.15"
$ws.Range("E45").Value = "  -7.35%  "

$ws.Range("E46").Value = "  -5.00%  "

$ws.Range("E47").Value = "  -7.16%  "

$ws.Range("E48").Value = "  -3.10%  "

$ws.Range("E49").Value = "  -4.24%  "

$ws.Range("E50").Value = "  -0.19%  "

$ws.Range("E51").Value = "  -6.13%  "
